$wb = $excel.ActiveWorkbook

$wsItems = $wb.Worksheets.Item("Items")
$wsEvents = $wb.Worksheets.Item("Events")

# --- Items sheet (sheet2.xml): update ItemID (column C) values for rows 2-9 ---
$wsItems.Range("C2").Value = 1001
$wsItems.Range("C3").Value = 1002
$wsItems.Range("C4").Value = 1003
$wsItems.Range("C5").Value = 1004
$wsItems.Range("C6").Value = 1005
$wsItems.Range("C7").Value = 1006
$wsItems.Range("C8").Value = 1007
$wsItems.Range("C9").Value = 1008

# --- Items sheet: append new row 10 (a new Dagger-like entry, ItemID 1009) ---
$wsItems.Range("A10").Value = 0
$wsItems.Range("B10").Value = 3020
$wsItems.Range("C10").Value = 1009
$wsItems.Range("D10").Value = "w"
$wsItems.Range("E10").Value = "p"
$wsItems.Range("F10").Value = "l"
$wsItems.Range("G10").Value = "c"
$wsItems.Range("H10").Value = 25
$wsItems.Range("I10").Value = 255
$wsItems.Range("J10").Value = 0
$wsItems.Range("K10").Value = 255
$wsItems.Range("L10").Value = "Dagger"
$wsItems.Range("M10").Value = 2
$wsItems.Range("N10").Value = 2
$wsItems.Range("O10").Value = 1.5
$wsItems.Range("P10").Value = 0
$wsItems.Range("Q10").Value = 0
$wsItems.Range("R10").Value = 0
$wsItems.Range("S10").Value = 0
$wsItems.Range("T10").Value = 0
$wsItems.Range("U10").Value = 0
$wsItems.Range("V10").Value = 0
$wsItems.Range("W10").Value = 0
$wsItems.Range("X10").Value = 0
$wsItems.Range("Y10").Value = 0
$wsItems.Range("Z10").Value = 0
$wsItems.Range("AA10").Value = 0
$wsItems.Range("AB10").Value = 1
$wsItems.Range("AC10").Value = 1
$wsItems.Range("AD10").Value = 1
$wsItems.Range("AE10").Value = 6
$wsItems.Range("AF10").Value = 1
$wsItems.Range("AG10").Value = 0
$wsItems.Range("AH10").Value = 0
$wsItems.Range("AI10").Value = 0
$wsItems.Range("AJ10").Value = 0
$wsItems.Range("AK10").Value = 0
$wsItems.Range("AL10").Value = 0
$wsItems.Range("AM10").Value = 0
$wsItems.Range("AN10").Value = 0
$wsItems.Range("AO10").Value = 0
$wsItems.Range("AP10").Value = 0
$wsItems.Range("AQ10").Value = 0
$wsItems.Range("AR10").Value = 0
$wsItems.Range("AS10").Value = 0
$wsItems.Range("AT10").Value = 0
$wsItems.Range("AU10").Value = 0
$wsItems.Range("AV10").Value = 0
$wsItems.Range("AW10").Value = 0
$wsItems.Range("AX10").Value = "A short knife with a pointed tip.&&(+1 Attack, +1 AC, +1 Dam, +1 MinDam, +1 maxDam, Pierce)"

# --- View/selection state ---
# Events sheet loses the active-tab flag; its selection moves to A2:J5.
$wsEvents.Activate()
$wsEvents.Range("A2:J5").Select()

# Items becomes the active sheet/tab, selection parked on the new row's C10 cell.
$wsItems.Activate()
$wsItems.Range("C10").Select()

Write-Output "done"
